$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the obsolete "上期/日盘" row (old row 3); remaining rows shift up.
$ws.Rows.Item(3).Delete()

# Update the remark for the 郑州/夜盘 row (now row 3) with the new note
# about strategy configuration files on jinrui servers.
$ws.Range("F3").Value = "其中`n(ok) fl34_cf_night, (ok)fl34_ta_night各给3手； `n(OK) fl50的各给1手`n，fl50的策略需要配置文件`n"

# Row grew taller to fit the new multi-line text.
$ws.Rows.Item(3).RowHeight = 104.25

# Move the active selection to match the new layout.
$ws.Range("E4").Select()
